$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 5: BEVNAT info
$ws.Cells.Item(5, 4).Value = "data/pdf/bevnat_info.pdf"
$ws.Cells.Item(5, 1).Value = "bevnat_info"
$ws.Cells.Item(5, 2).Value = "pdf"
$ws.Cells.Item(5, 3).Value = "BEVNAT: Fiche signalétique"
$ws.Cells.Item(5, 5).Value = 1706212962

# Row 6: STATPOP info
$ws.Cells.Item(6, 4).Value = "data/pdf/statpop_info.pdf"
$ws.Cells.Item(6, 3).Value = "STATPOP: Fiche signalétique"
$ws.Cells.Item(6, 1).Value = "statpop_info"
$ws.Cells.Item(6, 2).Value = "pdf"
$ws.Cells.Item(6, 5).Value = 1606212962

# Row 7: BEVNAT variable
$ws.Cells.Item(7, 4).Value = "data/pdf/bevnat_variable.pdf"
$ws.Cells.Item(7, 1).Value = "bevnat_variable"
$ws.Cells.Item(7, 2).Value = "pdf"
$ws.Cells.Item(7, 3).Value = "BEVNAT: Liste des variables"
$ws.Cells.Item(7, 5).Value = 1606212963

# Row 8: population press release
$ws.Cells.Item(8, 1).Value = "pop_com_1"
$ws.Cells.Item(8, 2).Value = "pdf"
$ws.Cells.Item(8, 4).Value = "data/pdf/pop_com_1.pdf"
$ws.Cells.Item(8, 3).Value = "Communiqué de presse population"
$ws.Cells.Item(8, 5).Value = 1724323867

# Widen column C to fit new, longer content
$ws.Columns.Item(3).ColumnWidth = 21.166666666666668

# Grow the table (ListObject) to cover the newly-added rows
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:E8"))

# Update active selection to reflect where editing left off
$ws.Range("E11").Select() | Out-Null
